# Trade #9 (internal trade index 38, "MarketMaking") closed at 2026-02-18 00:10:26.
# Two brand-new OPEN trades were also logged afterwards:
#   - trade 66 (MarketMaking) at 00:10:04
#   - trade 67 (momentum)      at 00:10:20
# This script reproduces every cell touched by that update across the
# Summary / Strategy Status / All Trades / momentum / MarketMaking sheets.

function Set-Text($ws, $row, $col, $text) {
    # Force text type so values that look like dates/times ("2026-02-18",
    # "00:10:04", ...) are not auto-coerced into date/time serials by Excel.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Set-Num($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - headline numbers after the close + the 2 new opens
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
Set-Num $wsSummary 3 2 1499.57   # Current Capital
Set-Num $wsSummary 4 2 0.67      # Total P&L $
Set-Num $wsSummary 5 2 0.36      # Total P&L %
Set-Num $wsSummary 6 2 37        # Total Trades
Set-Num $wsSummary 8 2 15        # Losing Trades
Set-Num $wsSummary 9 2 51.35     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
Set-Num $wsStatus 6 3 99.56999999999999   # Capital
Set-Num $wsStatus 6 4 8                   # Trades
Set-Num $wsStatus 6 5 -0.24               # P&L $
Set-Num $wsStatus 6 6 -0.43               # P&L %
Set-Num $wsStatus 6 7 37.5                # Win Rate %

# ---------------------------------------------------------------------
# 3) All Trades sheet
#    - row 39 (trade #38, MarketMaking) flips from OPEN -> CLOSED
#    - two new rows (67, 68) appended for the freshly opened trades
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

Set-Num  $wsAll 39 7  0.35                 # Exit Price
Set-Text $wsAll 39 8  "CLOSED"             # Status
Set-Num  $wsAll 39 9  -18.6047             # P&L %
Set-Num  $wsAll 39 10 -0.08                # P&L $
Set-Num  $wsAll 39 11 99.56999999999999    # Capital After
Set-Text $wsAll 39 12 "early_exit"         # Exit Reason
Set-Num  $wsAll 39 13 0.09                 # Duration (min)

# New row 67: trade 66, MarketMaking, OPEN
Set-Num  $wsAll 67 1  66
Set-Text $wsAll 67 2  "2026-02-18"
Set-Text $wsAll 67 3  "00:10:04"
Set-Text $wsAll 67 4  "MarketMaking"
Set-Text $wsAll 67 5  "UP"
Set-Num  $wsAll 67 6  0.47
Set-Text $wsAll 67 8  "OPEN"
Set-Num  $wsAll 67 9  0
Set-Num  $wsAll 67 10 0
Set-Num  $wsAll 67 11 99.65175839854133
Set-Num  $wsAll 67 13 0
Set-Num  $wsAll 67 14 0
Set-Num  $wsAll 67 15 0
Set-Num  $wsAll 67 16 0.6
Set-Text $wsAll 67 17 "Normal spread capture: 198 bps"

# New row 68: trade 67, momentum, OPEN
Set-Num  $wsAll 68 1  67
Set-Text $wsAll 68 2  "2026-02-18"
Set-Text $wsAll 68 3  "00:10:20"
Set-Text $wsAll 68 4  "momentum"
Set-Text $wsAll 68 5  "UP"
Set-Num  $wsAll 68 6  0.43
Set-Text $wsAll 68 8  "OPEN"
Set-Num  $wsAll 68 9  0
Set-Num  $wsAll 68 10 0
Set-Num  $wsAll 68 11 100
Set-Num  $wsAll 68 13 0
Set-Num  $wsAll 68 14 0
Set-Num  $wsAll 68 15 0
Set-Num  $wsAll 68 16 0.9
Set-Text $wsAll 68 17 "Upward momentum: 21.687% over 10 samples"

# ---------------------------------------------------------------------
# 4) momentum sheet - new row 5 for the freshly opened trade 67
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
Set-Num  $wsMomentum 5 1  67
Set-Text $wsMomentum 5 2  "2026-02-18"
Set-Text $wsMomentum 5 3  "00:10:20"
Set-Text $wsMomentum 5 4  "momentum"
Set-Text $wsMomentum 5 5  "UP"
Set-Num  $wsMomentum 5 6  0.43
Set-Text $wsMomentum 5 8  "OPEN"
Set-Num  $wsMomentum 5 9  0
Set-Num  $wsMomentum 5 10 0
Set-Num  $wsMomentum 5 11 100
Set-Num  $wsMomentum 5 12 0
Set-Num  $wsMomentum 5 13 0
Set-Num  $wsMomentum 5 14 0.9
Set-Text $wsMomentum 5 15 "Upward momentum: 21.687% over 10 samples"
Set-Num  $wsMomentum 5 17 0

# ---------------------------------------------------------------------
# 5) MarketMaking sheet
#    - row 10 (trade #38) flips from OPEN -> CLOSED
#    - new row 30 appended for the freshly opened trade 66
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

Set-Num  $wsMM 10 7  0.35
Set-Text $wsMM 10 8  "CLOSED"
Set-Num  $wsMM 10 9  -18.6047
Set-Num  $wsMM 10 10 -0.08
Set-Num  $wsMM 10 11 99.56999999999999
Set-Text $wsMM 10 16 "early_exit"
Set-Num  $wsMM 10 17 0.09

# New row 30: trade 66, MarketMaking, OPEN
Set-Num  $wsMM 30 1  66
Set-Text $wsMM 30 2  "2026-02-18"
Set-Text $wsMM 30 3  "00:10:04"
Set-Text $wsMM 30 4  "MarketMaking"
Set-Text $wsMM 30 5  "UP"
Set-Num  $wsMM 30 6  0.47
Set-Text $wsMM 30 8  "OPEN"
Set-Num  $wsMM 30 9  0
Set-Num  $wsMM 30 10 0
Set-Num  $wsMM 30 11 99.65175839854133
Set-Num  $wsMM 30 12 0
Set-Num  $wsMM 30 13 0
Set-Num  $wsMM 30 14 0.6
Set-Text $wsMM 30 15 "Normal spread capture: 198 bps"
Set-Num  $wsMM 30 17 0
